$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.379.57'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '1.788.76'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.22'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.556'
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.90'
$ws.Range('E8').Value = '  +3.45%  '
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0690'
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '2.046.09'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.14'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = '1.769.22'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.635'
$ws.Range('E15').Value = '  +2.07%  '
$ws.Range('D16').Value = '34.347.76'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.30'
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.49'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.33'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '0.0₃0795'
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('E21').Value = '  +3.21%  '
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '167.92'
$ws.Range('E24').Value = '  +3.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('E26').Value = '  +2.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.56'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.03'
$ws.Range('E30').Value = '  +8.31%  '
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('E32').Value = '  +2.47%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.62'
$ws.Range('E35').Value = '  +6.58%  '
$ws.Range('D36').Value = '1.408.45'
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.682'
$ws.Range('E37').Value = '  +4.84%  '
$ws.Range('E38').Value = '  +2.97%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '84.23'
$ws.Range('E40').Value = '  +4.46%  '
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.939'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.85'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0527'
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('E46').Value = '  +3.04%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').Value = '1.947.52'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.38'
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  -3.38%  '
